$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-31 12:56:21"

for ($row = 2; $row -le 529; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
